# Update cryptocurrency Price (D) and Volume(1h) (E) columns
# Values are written as literal text (matching the source XML's inlineStr
# cells) rather than being auto-converted by Excel into numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "290.19"
Set-TextValue $ws.Range("E2") "-4.04%"

Set-TextValue $ws.Range("D3") "31.57"
Set-TextValue $ws.Range("E3") "-1.05%"

Set-TextValue $ws.Range("D4") "4.954"
Set-TextValue $ws.Range("E4") "-3.14%"

Set-TextValue $ws.Range("D5") "0.07298"
Set-TextValue $ws.Range("E5") "-6.67%"

Set-TextValue $ws.Range("D6") "1.864"
Set-TextValue $ws.Range("E6") "-12.82%"

Set-TextValue $ws.Range("D7") "7.684"
Set-TextValue $ws.Range("E7") "-1.70%"

Set-TextValue $ws.Range("D8") "3.753"
Set-TextValue $ws.Range("E8") "-1.09%"

Set-TextValue $ws.Range("D9") "0.9109"
Set-TextValue $ws.Range("E9") "-1.13%"

Set-TextValue $ws.Range("D10") "0.1661"
Set-TextValue $ws.Range("E10") "-5.27%"

Set-TextValue $ws.Range("D11") "0.07680"
Set-TextValue $ws.Range("E11") "-0.50%"

Set-TextValue $ws.Range("D12") "0.08228"
Set-TextValue $ws.Range("E12") "-7.88%"

Set-TextValue $ws.Range("D13") "0.03024"
Set-TextValue $ws.Range("E13") "-2.92%"

Set-TextValue $ws.Range("D14") "0.1005"
Set-TextValue $ws.Range("E14") "0.24%"

Set-TextValue $ws.Range("D15") "0.001502"
Set-TextValue $ws.Range("E15") "-0.86%"

Set-TextValue $ws.Range("D16") "0.005918"
Set-TextValue $ws.Range("E16") "0.30%"

Set-TextValue $ws.Range("D18") "3.456"
Set-TextValue $ws.Range("E18") "0.19%"

Set-TextValue $ws.Range("D19") "2.130"
Set-TextValue $ws.Range("E19") "-6.62%"

Set-TextValue $ws.Range("D20") "0.3265"
Set-TextValue $ws.Range("E20") "-0.83%"

Set-TextValue $ws.Range("D21") "0.1296"
Set-TextValue $ws.Range("E21") "-2.51%"

Set-TextValue $ws.Range("D22") "4.375"
Set-TextValue $ws.Range("E22") "5.19%"

Set-TextValue $ws.Range("D23") "0.1999"
Set-TextValue $ws.Range("E23") "11.40%"

Set-TextValue $ws.Range("D24") "0.04496"
Set-TextValue $ws.Range("E24") "-2.04%"

Set-TextValue $ws.Range("D25") "0.001243"
Set-TextValue $ws.Range("E25") "0.36%"

Set-TextValue $ws.Range("D26") "0.004006"
Set-TextValue $ws.Range("E26") "-10.43%"

Set-TextValue $ws.Range("D27") "0.0001266"
Set-TextValue $ws.Range("E27") "1.17%"

Set-TextValue $ws.Range("D39") "0.01620"
Set-TextValue $ws.Range("E39") "-8.46%"

Set-TextValue $ws.Range("D40") "0.04404"
Set-TextValue $ws.Range("E40") "-7.60%"

Set-TextValue $ws.Range("D41") "0.007519"
Set-TextValue $ws.Range("E41") "5.15%"

Set-TextValue $ws.Range("D42") "0.1311"
Set-TextValue $ws.Range("E42") "-4.10%"

Set-TextValue $ws.Range("D43") "0.002359"
Set-TextValue $ws.Range("E43") "9.13%"

Set-TextValue $ws.Range("D44") "0.01096"
Set-TextValue $ws.Range("E44") "1.24%"

Set-TextValue $ws.Range("D45") "0.00006145"
Set-TextValue $ws.Range("E45") "-1.83%"

Set-TextValue $ws.Range("D46") "0.00000000760"
Set-TextValue $ws.Range("E46") "1.18%"

Set-TextValue $ws.Range("D47") "1.845"
Set-TextValue $ws.Range("E47") "57.27%"

Set-TextValue $ws.Range("D48") "0.003035"
Set-TextValue $ws.Range("E48") "-14.81%"

Set-TextValue $ws.Range("D49") "0.00002127"
Set-TextValue $ws.Range("E49") "1.18%"

Set-TextValue $ws.Range("D50") "0.0002025"
Set-TextValue $ws.Range("E50") "1.18%"

